# Non-Oncology Excel reports comparison
# Update the expected "ExcelReport" filename string for the Alkermes - Melanoma
# test case on the OldImportLogic sheet: remove the spaces around the dash
# before "Melanoma" ("ExcelReport-Alkermes - Melanoma-Economic-" ->
# "ExcelReport-Alkermes-Melanoma-Economic-").

$wb = $excel.ActiveWorkbook

$wsOld = $wb.Worksheets.Item("OldImportLogic")
$wsOld.Range("H3").Value = "ExcelReport-Alkermes-Melanoma-Economic-"

# Reflect the cell the user ended up editing/selecting as the active cell.
$wsOld.Activate()
$wsOld.Range("H3").Select()
